$d = $word.ActiveDocument
$W = "http://schemas.openxmlformats.org/wordprocessingml/2006/main"

# ------------------------------------------------------------------
# Locate the three paragraphs that need to change:
#   1) the outdated "Je suis inscrit pour le Windows ..." question
#   2) the paragraph right after it that only carries the _GoBack
#      bookmark
#   3) the empty paragraph that follows the bookmark paragraph
# We find paragraph (1) by its text instead of a hard-coded index so
# the script keeps working even if earlier content in the document
# shifts paragraph numbers around.
# ------------------------------------------------------------------
$finder = $d.Content
$found = $finder.Find.Execute(
    "Je suis inscrit pour le Windows", $true, $false, $false, $false,
    $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find the 'Je suis inscrit pour le Windows ...' paragraph"
}

$questionPara  = $finder.Paragraphs(1).Range
$bookmarkPara  = $questionPara.Next(4, 1)
$trailingPara  = $bookmarkPara.Next(4, 1)

# Range covering all three paragraphs, from the start of the question
# through the end of the trailing empty paragraph.
$block = $d.Range($questionPara.Start, $trailingPara.End)

# ------------------------------------------------------------------
# Replace that block with:
#   - one paragraph with a bold "NB : " lead-in followed by the
#     English note about the Singleton pattern (English language
#     formatting applied throughout, including the paragraph mark),
#     keeping the existing _GoBack bookmark on that same paragraph
#   - two empty English-tagged paragraphs where there used to be one
# ------------------------------------------------------------------
$xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
       '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
       '<pkg:xmlData>' +
       '<w:document xmlns:w="' + $W + '"><w:body>' +
           '<w:p>' +
               '<w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr>' +
               '<w:r><w:rPr><w:b/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">NB : </w:t></w:r>' +
               '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>The Singleton design pat</w:t></w:r>' +
               '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>tern is used on the connection string.</w:t></w:r>' +
               '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>' +
           '</w:p>' +
           '<w:p><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr></w:p>' +
           '<w:p><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr></w:p>' +
       '</w:body></w:document>' +
       '</pkg:xmlData></pkg:part></pkg:package>'

$block.InsertXML($xml)
